$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the hyperlink row (currently row 9) so the
# hyperlink row becomes row 10.
$ws.Rows.Item(9).Insert()

# Fill in new sample data for the inserted row 9
$ws.Range("A9").Value = 45
$ws.Range("B9").Value = 180

# Copy the formulas from row 8 down into the new row 9, the same
# formulas used throughout column C:J (just like dragging the fill
# handle down one more row in Excel).
$ws.Range("C9").Formula = "=SIN(RADIANS(A9))"
$ws.Range("D9").Formula = "=COS(RADIANS(A9))"
$ws.Range("E9").Formula = "=SIN(RADIANS(B9))"
$ws.Range("F9").Formula = "=COS(RADIANS(B9))"
$ws.Range("G9").Formula = "=+C9*F9-D9*E9"
$ws.Range("H9").Formula = "=+C9*E9+D9*F9"
$ws.Range("I9").Formula = '=IF(G9<0, "CW", "CCW")'
$ws.Range("J9").Formula = "=-DEGREES(ATAN2(H9,G9))"

# Remove the broken array formula that used to live in C9 and now,
# after the row insert, has shifted down to C10.
$ws.Range("C10").ClearContents()

# Fix up the hyperlink so it references the shifted hyperlink row (10)
# instead of the old row 9.
$ws.Range("A9").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A10"), "https://stackoverflow.com/a/16544330/14775294") | Out-Null
$ws.Range("A10").Style = "Hyperlink"

# Update the selection to match the new active cell
$ws.Range("E9").Select()
